# Applies the "Development Process" slide content edit:
# Slide 10, shape "Content Placeholder 2" gets populated with the
# project retrospective text (previously just an empty paragraph).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(10)

$sh = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $cand = $s.Shapes.Item($i)
    if ($cand.Name -eq "Content Placeholder 2") { $sh = $cand }
}
if ($sh -eq $null) { $sh = $s.Shapes.Item(2) }

$tf = $sh.TextFrame
$tr = $tf.TextRange

$p1r1 = "Requirement "
$p1r2 = "analysis->program design->implementation->AI integration ->Function improvements->"
$p1r3 = "verify"

$p2r1 = "Our project idea came from several programming assignments on planning travel routes this semester. we want to create a travel planner closer to real-life travel. "

$p3r1 = "Based "
$p3r2 = "on the second job, anytime search, we extended the algorithm to enable it to plan daily itineraries and introduced "
$p3r3 = "neural network prediction to "
$p3r4 = "select routes based on the user's interests"
$p3r5 = "."

$p4r1 = "ChatGPT4.0 was used to assist in the design of neural network. There were many problems with the "
$p4r2 = "initial "
$p4r3 = "design, but after constantly clarifying my needs, the AI worked "
$p4r4 = "satisfactorily."

$p5r1 = "During the development of this project, we became more familiar with search algorithms and neural network training. In future expectations, we may extend the front-end of this project to implement a graphical interface for user "
$p5r2 = "interaction."

$para2 = $p2r1
$para3 = $p3r1 + $p3r2 + $p3r3 + $p3r4 + $p3r5
$para4 = $p4r1 + $p4r2 + $p4r3 + $p4r4
$para5 = $p5r1 + $p5r2

$tr.Text = $p1r1 + "`r" + $para2 + "`r" + $para3 + "`r" + $para4 + "`r" + $para5

# Paragraph 1: plain (no bullet), left margin/indent reset to 0. Build it
# out of three separate runs via InsertAfter (mirrors the three <a:r>
# runs in the source edit, which differ only by transient autocorrect
# state we can't reproduce, not by visible formatting).
$par1 = $tr.Paragraphs(1, 1)
$par1.ParagraphFormat.Bullet.Visible = 0
$rl = $tf.Ruler.Levels(1)
$rl.LeftMargin = 0
$rl.FirstMargin = 0
$par1.InsertAfter($p1r2) | Out-Null
$par1.InsertAfter($p1r3) | Out-Null

# Paragraphs 2-5 use a smaller 20pt font size; set per-run to mirror the
# multi-run structure of the source edit (several runs per paragraph).
$par2 = $tr.Paragraphs(2, 1)
$par2.Characters(1, $p2r1.Length).Font.Size = 20

$par3 = $tr.Paragraphs(3, 1)
$off = 1
$par3.Characters($off, $p3r1.Length).Font.Size = 20
$off += $p3r1.Length
$par3.Characters($off, $p3r2.Length).Font.Size = 20
$off += $p3r2.Length
$par3.Characters($off, $p3r3.Length).Font.Size = 20
$off += $p3r3.Length
$par3.Characters($off, $p3r4.Length).Font.Size = 20
$off += $p3r4.Length
$par3.Characters($off, $p3r5.Length).Font.Size = 20

$par4 = $tr.Paragraphs(4, 1)
$off = 1
$par4.Characters($off, $p4r1.Length).Font.Size = 20
$off += $p4r1.Length
$par4.Characters($off, $p4r2.Length).Font.Size = 20
$off += $p4r2.Length
$par4.Characters($off, $p4r3.Length).Font.Size = 20
$off += $p4r3.Length
$par4.Characters($off, $p4r4.Length).Font.Size = 20

$par5 = $tr.Paragraphs(5, 1)
$off = 1
$par5.Characters($off, $p5r1.Length).Font.Size = 20
$off += $p5r1.Length
$par5.Characters($off, $p5r2.Length).Font.Size = 20

# The shape now autofits its text (shrink text on overflow), matching the
# <a:normAutofit/> seen in the edited deck.
$tf.AutoSize = 2
